$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 163; existing rows 163-168 shift down to 164-169,
# carrying their formatting (including the date style on column D) with them.
$ws.Rows.Item(163).Insert()

# Populate the newly inserted row 163 with the new weekly record.
$ws.Range("A163").Value = 5
$ws.Range("B163").Value = "Macroferia Regional de Talca"
$ws.Range("C163").Value = "Maule"
$ws.Range("D163").Value = 45041
$ws.Range("E163").Value = 7
$ws.Range("F163").Value = 100112001
$ws.Range("G163").Value = "Berenjena"
$ws.Range("H163").Value = "Sin especificar"
$ws.Range("I163").Value = "Primera"
$ws.Range("J163").Value = 150
$ws.Range("K163").Value = 10000
$ws.Range("L163").Value = 10000
$ws.Range("M163").Value = 10000
$ws.Range("N163").Value = "$/caja 50 unidades"
$ws.Range("O163").Value = "Región del Maule"
$ws.Range("P163").Value = 200
$ws.Range("Q163").Value = 50
$ws.Range("R163").Value = "Hortaliza"
